$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text formatting on Price/Volume columns so Excel does not
# auto-convert numeric-looking strings (stripping trailing zeros, etc.)
$ws.Range("D2:E51").NumberFormat = "@"

$updates = @{
    "D2" = "69.231.87"
    "E2" = "  -2.75%  "
    "D3" = "3.677.12"
    "E3" = "  -3.21%  "
    "E4" = "  +0.06%  "
    "D5" = "680.40"
    "E5" = "  -3.74%  "
    "D6" = "159.22"
    "E6" = "  -6.89%  "
    "D7" = "3.675.41"
    "E7" = "  -3.22%  "
    "E8" = "  -0.03%  "
    "D9" = "0.493"
    "E9" = "  -6.16%  "
    "D10" = "0.144"
    "E10" = "  -9.81%  "
    "D11" = "7.08"
    "E11" = "  -5.63%  "
    "D12" = "0.434"
    "E12" = "  -9.43%  "
    "D13" = "0.0000231"
    "E13" = "  -7.67%  "
    "D14" = "4.298.65"
    "E14" = "  -3.18%  "
    "D15" = "32.28"
    "E15" = "  -10.76%  "
    "D16" = "3.670.11"
    "E16" = "  -2.81%  "
    "D17" = "69.285.35"
    "E17" = "  -2.74%  "
    "E18" = "  -1.48%  "
    "D19" = "15.78"
    "E19" = "  -9.73%  "
    "D20" = "6.40"
    "E20" = "  -10.66%  "
    "D21" = "470.50"
    "E21" = "  -8.79%  "
    "D22" = "9.83"
    "E22" = "  -5.53%  "
    "D23" = "0.646"
    "E23" = "  -9.46%  "
    "D24" = "79.20"
    "E24" = "  -5.62%  "
    "D25" = "3.824.25"
    "E25" = "  -3.00%  "
    "E26" = "  +0.01%  "
    "D27" = "0.0000123"
    "E27" = "  -11.96%  "
    "D28" = "10.83"
    "E28" = "  -14.49%  "
    "D29" = "9.05"
    "E29" = "  -12.35%  "
    "D30" = "2.68"
    "E30" = "  -11.37%  "
    "D31" = "1.72"
    "E31" = "  -14.75%  "
    "D32" = "6.64"
    "E32" = "  -9.83%  "
    "D33" = "2.01"
    "E33" = "  -10.18%  "
    "E34" = "  +0.38%  "
    "D35" = "26.54"
    "E35" = "  -8.80%  "
    "B36" = "Kaspa"
    "C36" = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
    "D36" = "0.160"
    "E36" = "  -7.93%  "
    "B37" = "Aptos"
    "C37" = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
    "D37" = "8.08"
    "E37" = "  -12.54%  "
    "B38" = "Filecoin"
    "C38" = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
    "D38" = "6.03"
    "E38" = "  -6.77%  "
    "B39" = "Stacks"
    "C39" = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
    "D39" = "2.24"
    "E39" = "  -7.98%  "
    "B40" = "USDe"
    "C40" = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
    "D40" = "1.00"
    "E40" = "  -0.01%  "
    "B41" = "Hedera"
    "C41" = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
    "D41" = "0.0898"
    "E41" = "  -11.10%  "
    "B42" = "FirstDigitalUSD"
    "C42" = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
    "D42" = "1.00"
    "E42" = "  +0.01%  "
    "B43" = "Mantle"
    "C43" = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
    "D43" = "0.939"
    "E43" = "  -7.01%  "
    "B44" = "Monero"
    "C44" = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
    "D44" = "165.70"
    "E44" = "  -1.66%  "
    "B45" = "OKB"
    "C45" = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
    "D45" = "47.68"
    "E45" = "  -4.77%  "
    "B46" = "dogwifhat"
    "C46" = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
    "D46" = "2.69"
    "E46" = "  -17.02%  "
    "B47" = "ONDO"
    "C47" = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
    "D47" = "1.29"
    "E47" = "  -7.61%  "
    "D48" = "27.91"
    "E48" = "  -3.64%  "
    "B49" = "SuiNetwork"
    "C49" = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
    "D49" = "1.08"
    "E49" = "  -6.21%  "
    "B50" = "FLOKI"
    "C50" = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
    "D50" = "0.000268"
    "E50" = "  -12.49%  "
    "B51" = "Cosmos"
    "C51" = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
    "D51" = "7.85"
    "E51" = "  -8.66%  "
}

foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}
